# Adds Balance/Date/Сurrency/Year/Month columns (E:I) to the bank-operations
# table, converts the Card column for the "unique card" summary rows (8-11)
# to plain text, and appends an extra duplicate row (15) after the existing
# totals rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (E1:I1) - same bordered/centred style as B1:D1.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 5).Value = "Balance"
$ws.Cells.Item(1, 6).Value = "Date"
$ws.Cells.Item(1, 7).Value = "Сurrency"
$ws.Cells.Item(1, 8).Value = "Year"
$ws.Cells.Item(1, 9).Value = "Month"

$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Range("E1:I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Detail rows 2-7: add Balance / Date / Currency / Year / Month, and make
# the Card column (C) a real number (it used to be a shared string).
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 3).Value = 6678
$ws.Cells.Item(2, 5).Value = 588
$ws.Cells.Item(2, 6).Value = 43375.62273148148
$ws.Cells.Item(2, 7).Value = "EUR"
$ws.Cells.Item(2, 8).Value = 2018
$ws.Cells.Item(2, 9).Value = 10

$ws.Cells.Item(3, 3).Value = 1238
$ws.Cells.Item(3, 5).Value = 450
$ws.Cells.Item(3, 6).Value = 43409.44231481481
$ws.Cells.Item(3, 7).Value = "EUR"
$ws.Cells.Item(3, 8).Value = 2018
$ws.Cells.Item(3, 9).Value = 11

$ws.Cells.Item(4, 3).Value = 6678
$ws.Cells.Item(4, 5).Value = 578
$ws.Cells.Item(4, 6).Value = 43436.62273148148
$ws.Cells.Item(4, 7).Value = "EUR"
$ws.Cells.Item(4, 8).Value = 2018
$ws.Cells.Item(4, 9).Value = 12

$ws.Cells.Item(5, 3).Value = 6623
$ws.Cells.Item(5, 5).Value = 870
$ws.Cells.Item(5, 6).Value = 43437.41219907408
$ws.Cells.Item(5, 7).Value = "EUR"
$ws.Cells.Item(5, 8).Value = 2018
$ws.Cells.Item(5, 9).Value = 12

$ws.Cells.Item(6, 3).Value = 6678
$ws.Cells.Item(6, 5).Value = 548
$ws.Cells.Item(6, 6).Value = 43438.62273148148
$ws.Cells.Item(6, 7).Value = "EUR"
$ws.Cells.Item(6, 8).Value = 2018
$ws.Cells.Item(6, 9).Value = 12

$ws.Cells.Item(7, 3).Value = 1253
$ws.Cells.Item(7, 5).Value = 700
$ws.Cells.Item(7, 6).Value = 43467.75717592592
$ws.Cells.Item(7, 7).Value = "EUR"
$ws.Cells.Item(7, 8).Value = 2019
$ws.Cells.Item(7, 9).Value = 1

$ws.Range("F2:F7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Rows 8-11: "unique card" summary rows. Card column keeps its original
# text type (e.g. "6678"), only Balance + Currency are new. A leading
# apostrophe forces Excel to store the digits as text instead of a
# number; resetting the style afterwards drops the quote-prefix flag
# that the apostrophe entry leaves behind so the cell ends up identical
# to a plain shared-string cell (no explicit style index).
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 3).Value = "'6678"
$ws.Cells.Item(8, 5).Value = 548
$ws.Cells.Item(8, 7).Value = "EUR"

$ws.Cells.Item(9, 3).Value = "'6623"
$ws.Cells.Item(9, 5).Value = 870
$ws.Cells.Item(9, 7).Value = "EUR"

$ws.Cells.Item(10, 3).Value = "'1238"
$ws.Cells.Item(10, 5).Value = 450
$ws.Cells.Item(10, 7).Value = "EUR"

$ws.Cells.Item(11, 3).Value = "'1253"
$ws.Cells.Item(11, 5).Value = 700
$ws.Cells.Item(11, 7).Value = "EUR"

$ws.Range("C8:C11").Style = "Normal"

# ---------------------------------------------------------------------
# Row 12: Total row - Balance + Currency only (Bank/Card/Operation blank)
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 5).Value = 2568
$ws.Cells.Item(12, 7).Value = "EUR"

# ---------------------------------------------------------------------
# Rows 13-14: the leading index column goes back to a plain 0/1 number
# (it used to hold the "Total"/"GorgeousBank" shared strings) and the
# Card column becomes numeric again.
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = 0
$ws.Cells.Item(13, 2).Value = "SuperBank"
$ws.Cells.Item(13, 3).Value = 6623

$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "SuperBank"
$ws.Cells.Item(14, 3).Value = 6678

# ---------------------------------------------------------------------
# Row 15: brand-new row, duplicate of row 14.
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = 0
$ws.Cells.Item(15, 2).Value = "SuperBank"
$ws.Cells.Item(15, 3).Value = 6678
$ws.Cells.Item(15, 4).Value = -40

# Column-A cells carry the bordered/centred header style (s="1"); copy it
# from the row above instead of guessing the Excel style name.
$ws.Cells.Item(14, 1).Copy() | Out-Null
$ws.Cells.Item(15, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
